$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The "Requisitos" list paragraph (the ListBullet paragraph right
# after the "Requisitos" heading) is rewritten: the course-requirement
# lines are reordered, a few lines are dropped, a couple of codes /
# diacritics are corrected, and the remaining lines keep the original
# "<code> -  <description>  (Requisito)" text followed by a manual
# line break (w:br), one run per line - matching how the paragraph was
# originally authored.
#
# Approach: type the full replacement list as brand-new runs at the
# very end of the document (the Requisitos list is the document's
# last paragraph, so Selection.EndKey(wdStory) lands right before its
# paragraph mark). Typing each line as its own TypeText call keeps
# every line in its own run, exactly like the source markup. Then
# remove every original line (its text plus trailing line break) via
# Find/Replace (wdReplaceOne) so only the old runs disappear - the
# freshly typed runs are left completely untouched even when a line's
# text is unchanged and therefore appears twice in the document while
# both copies are present.
# ------------------------------------------------------------------

$sel = $word.Selection
$sel.EndKey(6)

$sel.TypeText('LOM3105 -  Computação e análise de dados em Engenharia  (Requisito)' + [char]11)
$sel.TypeText('LOQ4095 -  Química Geral Experimental  (Requisito)' + [char]11)
$sel.TypeText('LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)' + [char]11)
$sel.TypeText('LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)' + [char]11)
$sel.TypeText('LOB1046 -  Engenharia do Meio Ambiente  (Requisito)' + [char]11)
$sel.TypeText('LOB1053 -  Física III  (Requisito)' + [char]11)
$sel.TypeText('LOB1003 -  Cálculo I  (Requisito)' + [char]11)
$sel.TypeText('LOB1012 -  Estatística  (Requisito)' + [char]11)
$sel.TypeText('LOB1036 -  Geometria Analítica  (Requisito)' + [char]11)
$sel.TypeText('LOB1037 -  Álgebra Linear  (Requisito)' + [char]11)
$sel.TypeText('LOB1038 -  Física Experimental I  (Requisito)' + [char]11)
$sel.TypeText('LOB1039 -  Física Experimental III  (Requisito)' + [char]11)
$sel.TypeText('LOB1041 -  Física Experimental II  (Requisito)' + [char]11)
$sel.TypeText('LOB1052 -  Cálculo III  (Requisito)' + [char]11)
$sel.TypeText('LOM3037 -  Química Inorgânica  (Requisito)' + [char]11)
$sel.TypeText('LOM3056 -  Fundamentos de Química Orgânica  (Requisito)' + [char]11)
$sel.TypeText('LOM3099 -  Estática  (Requisito)' + [char]11)
$sel.TypeText('LOQ4264 -  Engenharia da Sustentabilidade  (Requisito)' + [char]11)
$sel.TypeText('LOB1004 -  Cálculo II  (Requisito)' + [char]11)
$sel.TypeText('LOB1018 -  Física I  (Requisito)' + [char]11)
$sel.TypeText('LOB1019 -  Física II  (Requisito)' + [char]11)
$sel.TypeText('LOM3018 -  Introdução à Engenharia de Materiais  (Requisito)' + [char]11)
$sel.TypeText('LOM3013 -  Ciência dos Materiais  (Requisito)' + [char]11)
$sel.TypeText('LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)' + [char]11)

$rng = $d.Content

$rng.Find.Execute('LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1053 -  Física III  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1038 -  Física Experimental I  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3013 -  Ciência dos Materiais  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1018 -  Física I  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1041 -  Física Experimental II  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOQ4095 -  Química Geral Experimental  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1036 -  Geometria Analítica  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3037 -  Química Inorgânica  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3099 -  Estática  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3056 -  Fundamentos de Química Orgânica  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOQ4246 -  Engenharia da Sustentabilidade  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3018 -  Introdução à Engenharia de Materiais  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1012 -  Estatística  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1004 -  Cálculo II  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1046 -  Engenharia do Meio Ambiente  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1003 -  Cálculo I  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1052 -  Cálculo III  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1037 -  Àlgebra Linear  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1019 -  Física II  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOB1039 -  Física Experimental III  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
$rng.Find.Execute('LOM3105 -  Computação e análise de dados em Engenharia  (Requisito)' + [char]11, $true, $false, $false, $false, $false, $true, 1, $false, '', 1) | Out-Null
